# VERSION DE PRUEBAS NO A LA VISTA DEL CLIENTE
#
# 1) "Vamos a hacer pruebas con git" -> split so "git" is wrapped in
#    proofErr spellStart/spellEnd (keeps the existing "…." run intact).
# 2) "...de un modulo separado..." -> split so "modulo" is wrapped in
#    proofErr spellStart/spellEnd.
# 3) Insert a new underlined paragraph after "Versión de explotación"
#    with the bookmark ( _GoBack ) moved onto it, and drop the trailing
#    empty paragraph.

$d = $word.ActiveDocument
$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

# --- helper: locate the paragraph whose text starts with $prefix -----
function Get-ParaByPrefix($prefix) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $t = $d.Paragraphs($i).Range.Text
        if ($t.StartsWith($prefix)) {
            return $d.Paragraphs($i)
        }
    }
    return $null
}

# --- 1) "Vamos a hacer pruebas con git" ------------------------------
$p1 = Get-ParaByPrefix("Vamos a hacer pruebas con git")
$r1 = $p1.Range
$xml1 = @"
<w:p xmlns:w='$wNs'>
<w:r><w:t xml:space='preserve'>Vamos a hacer pruebas con </w:t></w:r>
<w:proofErr w:type='spellStart'/>
<w:r><w:t>git</w:t></w:r>
<w:proofErr w:type='spellEnd'/>
<w:r><w:t>….</w:t></w:r>
</w:p>
"@
$r1.InsertXML($xml1) | Out-Null

# --- 2) "...de un modulo separado del tronco principal" --------------
$p2 = Get-ParaByPrefix("Voy a desarrollar por mi cuenta")
$r2 = $p2.Range
$xml2 = @"
<w:p xmlns:w='$wNs'>
<w:r><w:t xml:space='preserve'>Voy a desarrollar por mi cuenta otra rama para probar la funcionalidad de un </w:t></w:r>
<w:proofErr w:type='spellStart'/>
<w:r><w:t>modulo</w:t></w:r>
<w:proofErr w:type='spellEnd'/>
<w:r><w:t xml:space='preserve'> separado del tronco principal</w:t></w:r>
</w:p>
"@
$r2.InsertXML($xml2) | Out-Null

# --- 3) New hidden-from-client paragraph + bookmark move --------------
$p3 = Get-ParaByPrefix("Versión de explotación")
$p3Next = $p3.Next()
if ($p3Next -ne $null -and $p3Next.Range.Text.Trim() -eq "") {
    # the trailing empty paragraph right after it is being dropped
    $endPara = $p3Next
} else {
    $endPara = $p3
}

$start = $p3.Range.Start
$end = $endPara.Range.End
$r3 = $d.Range($start, $end)
$xml3 = @"
<w:p xmlns:w='$wNs' w:rsidR='00894EAC' w:rsidRDefault='00894EAC'>
<w:pPr>
<w:pBdr>
<w:top w:val='single' w:sz='6' w:space='1' w:color='auto'/>
<w:bottom w:val='single' w:sz='6' w:space='1' w:color='auto'/>
</w:pBdr>
</w:pPr>
<w:r><w:t>Versión de explotación</w:t></w:r>
</w:p>
<w:p xmlns:w='$wNs'>
<w:pPr>
<w:rPr>
<w:u w:val='single'/>
</w:rPr>
</w:pPr>
<w:r><w:t>Este es mi nuevo párrafo que no quiero que el cliente no vea</w:t></w:r>
<w:bookmarkStart w:id='0' w:name='_GoBack'/>
<w:bookmarkEnd w:id='0'/>
</w:p>
"@
$r3.InsertXML($xml3) | Out-Null
